$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns before column B, shifting the existing B:J data (and
# everything below) over to K:S.
$ws.Range("B1:J1").EntireColumn.Insert()

# Keep the newly inserted columns (and the date columns they pushed to the
# right) the same width (8 characters) as the rest of the sheet. Column B
# (like before the edit) is left at the default width.
$ws.Range("C1:S1").EntireColumn.ColumnWidth = 7.1666666666667

# New header labels for the 9 freshly inserted columns (newest-to-oldest,
# continuing the existing left-to-right newest-first layout).
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# Fill the new data columns (rows 2-33) with the same "UN" placeholder used
# throughout the rest of the sheet.
$ws.Range("B2:J33").Value = "UN"
